$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03435667493040919
$ws.Range("C2").Value = 0.9996458930772465
$ws.Range("D2").Value = 0.1335970946467847
$ws.Range("G2").Value = 0.1217136106832186
$ws.Range("H2").Value = 0.9740000000000001
